$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Lab of Ornithology Reference (Non-Circulating)"
$ws.Range("A3").Value = "Lab of Ornithology"
$ws.Range("E3").Value = "Lab of Ornithology > Main Collection"
$ws.Range("E4").Value = "Lab of Ornithology > Reference"
$ws.Range("E5").Value = "Lab of Ornithology > Museum of Vertebrates"
$ws.Range("E6").Value = "Lab of Ornithology > Macaulay Library"

$ws.Range("D6").Select() | Out-Null
